# Add new participant "240M_FM" as the next row in the ID/Order list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells(1, 1).End(-4121).Row  # xlDown
$newRow = $lastRow + 1
$newOrder = $ws.Cells($lastRow, 2).Value2 + 1

$ws.Cells($newRow, 1).Value = "240M_FM"
$ws.Cells($newRow, 2).Value = $newOrder

# Reflect the printable-page setup Excel stamps onto the sheet once a
# user interacts with Page Setup / Print.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Scroll/selection state left behind after entering the new row.
$excel.ActiveWindow.ScrollRow = 70
$ws.Range("G86").Select()
